$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cranking")

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P")

# --- Row 4: fuel enrichment is now capped at 130 (IF(...<130,...,130)) ---
# B4 holds its own (non-shared) formula; C4 is the anchor of the shared
# formula group C4:P4 - writing the whole range at once keeps the group
# (t="shared") intact instead of exploding it into per-cell formulas.
$ws.Range("B4").Formula = "=IF(B2+B2*B3<130,B2+B2*B3,130)"
$ws.Range("C4:P4").Formula = "=IF(C2+C2*C3<130,C2+C2*C3,130)"

# --- New row 5: "Fine tuning Addition" (same look as row 3's Addition row) ---
$ws.Range("A5").Value = "Fine tuning Addition"
$ws.Range("B3:P3").Copy() | Out-Null
$ws.Range("B5:P5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$row5vals = @(5, 5, 5, 2.5, 1, 0.5, 0.25, 0.125, 0, 0, 0, 0, 0, 0, 0)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "5").Value = $row5vals[$i]
}

# --- New row 6: "E85 Cranking WIP" - capped enrichment using row 5 addition ---
$ws.Range("A6").Value = "E85 Cranking WIP"
$ws.Range("B4:P4").Copy() | Out-Null
$ws.Range("B6:P6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

foreach ($c in $cols) {
    $ws.Range($c + "6").Formula = "=IF(" + $c + "2+" + $c + "2*" + $c + "5<130," + $c + "2+" + $c + "2*" + $c + "5,130)"
}

# Column A is now wide enough to show the new row labels in full.
$ws.Columns.Item(1).AutoFit() | Out-Null

# --- View state: Cranking becomes the active/selected tab, cursor at A7 ---
$ws.Activate()
$ws.Range("A7").Select() | Out-Null

$win = $excel.ActiveWindow
$win.Left = -57720
